# Sprint tracker update: refresh sprint data (sprints 7/8) and
# "Monografia" edits, drop the stray empty "Coluna1" table column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint")

# --- New table content -----------------------------------------------
# Row 1 = headers (unchanged); rows 2-4 = refreshed task log.
$data = @(
    @("DATA", "TAREFA REALIZADA ", "QUEM REALIZOU", "STATUS"),
    @("17/04/2024 - 24/04/2024", "Edições na Tela de Fornecedores", "Gabriel e Cristielen", "Pronto"),
    @("17/04/2024 - 24/04/2024", "Edições no Manual do Usuário", "Gabriel", "Pronto"),
    @("17/04/2024 - 24/04/2024", "Edições na Monografia", "Bruno", "Pronto")
)

for ($r = 0; $r -lt $data.Length; $r++) {
    for ($c = 0; $c -lt $data[$r].Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $data[$r][$c]
    }
}

# Drop the old, unused "Coluna1" column's leftover cell content (column E).
$ws.Range("E1:E6").Clear()

# Re-apply the plain centered style (same look as the header row) across
# the whole refreshed block so every cell is visually consistent again.
$ws.Range("A1").Copy()
$ws.Range("A1:D4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Shrink the table to match the new 4-column x 4-row extent.
$lo = $ws.ListObjects.Item(1)
$lo.ListColumns.Item(5).Delete()
$lo.Resize($ws.Range("A1:D4"))

# Matches the saved selection in the edited workbook.
$null = $ws.Range("C13").Select()
